$wb = $excel.ActiveWorkbook

# --- Settings sheet ---
$wsSettings = $wb.Worksheets.Item("Settings")
# Rename the queue referenced by OrchestratorQueueName
$wsSettings.Range("B2").Value = "CandidatesQueue"
# Fill in the previously-empty OrchestratorQueueFolder value
$wsSettings.Range("B3").Value = "OnboardingEmailCoordinator"
$wsSettings.Range("B3").Select() | Out-Null

# --- Constants sheet ---
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("A19").Value = "EmailSubject"
$wsConstants.Range("B19").Value = "Automation Error!"
$wsConstants.Range("A20").Value = "EmailBody"
$wsConstants.Range("B20").Value = "Hello, `nAn exception occurred during the automation process.  Please find the details below:`nException Source: @Source`nException Message: @Message`nA screenshot of the error has been attached for reference. Please see the attachment for more details.`nThank you and have a good day,`nRobot :)"
$wsConstants.Range("B20").WrapText = $true
$wsConstants.Rows.Item(20).RowHeight = 14.25
$wsConstants.Range("A20").Select() | Out-Null

# --- Assets sheet ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("A2").Value = "EmailCredentials"
$wsAssets.Range("B2").Value = "EmailCredentials"

$wsAssets.Range("A3").Value = "GetAverageProcessingTimeQuery"
$wsAssets.Range("B3").Value = "GetAverageProcessingTimeQuery"

$wsAssets.Range("A4").Value = "GetLongRunningTransactionsQuery"
$wsAssets.Range("B4").Value = "GetLongRunningTransactionsQuery"

$wsAssets.Range("A5").Value = "GetTransactionStatusAggregatesQuery"
$wsAssets.Range("B5").Value = "GetTransactionStatusAggregatesQuery"

$wsAssets.Range("A6").Value = "GetDateRangeQuery"
$wsAssets.Range("B6").Value = "GetDateRangeQuery"

$wsAssets.Range("B6").Select() | Out-Null
